# Atualização da base de dados
# Update "Inscritos" (column E) counts for specific rows per the data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E6").Value  = 23
$ws.Range("E8").Value  = 22
$ws.Range("E15").Value = 47
$ws.Range("E16").Value = 171
$ws.Range("E18").Value = 38
